$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4065.25
$ws.Range("I98").Value = 3674.75
$ws.Range("J98").Value = 4846.25
$ws.Range("K98").Value = 3674.75
$ws.Range("L98").Value = 4846.25
$ws.Range("M98").Value = -2176.75
$ws.Range("N98").Value = -7842.25
$ws.Range("H113").Value = 3368.9
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 3543.2222
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 3543.2222
$ws.Range("M113").Value = 1454
$ws.Range("N113").Value = -10051.2222
$ws.Range("H122").Value = 4065.25
$ws.Range("I122").Value = 3674.75
$ws.Range("J122").Value = 4846.25
$ws.Range("K122").Value = 11024.25
$ws.Range("L122").Value = 14538.75
$ws.Range("M122").Value = -8574.25
$ws.Range("N122").Value = -19438.75
$ws.Range("H138").Value = 1346.8096
$ws.Range("I138").Value = 946.69696
$ws.Range("J138").Value = 1786.9333
$ws.Range("K138").Value = 2840.09088
$ws.Range("L138").Value = 5360.7999
$ws.Range("M138").Value = 2299.90912
$ws.Range("N138").Value = -15640.7999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1263.5333
$ws.Range("I61").Value = 957.38464
$ws.Range("J61").Value = 3253.5
$ws.Range("K61").Value = 957.38464
$ws.Range("L61").Value = 3253.5
$ws.Range("M61").Value = -745.38464
$ws.Range("N61").Value = -3677.5
$ws.Range("H74").Value = 1037.619
$ws.Range("I74").Value = 820.8889
$ws.Range("K74").Value = 820.8889
$ws.Range("M74").Value = 53.11109999999996
$ws.Range("H77").Value = 1037.619
$ws.Range("I77").Value = 820.8889
$ws.Range("K77").Value = 4104.444500000001
$ws.Range("M77").Value = 263.5554999999995
$ws.Range("H102").Value = 15152980
$ws.Range("I102").Value = 15152980
$ws.Range("K102").Value = 15152980
$ws.Range("M102").Value = -15151358
$ws.Range("H122").Value = 2059.1
$ws.Range("I122").Value = 2059.1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6177.299999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3727.299999999999
$ws.Range("N122").Value = $null
$ws.Range("H136").Value = 1263.5333
$ws.Range("I136").Value = 957.38464
$ws.Range("J136").Value = 3253.5
$ws.Range("K136").Value = 2872.15392
$ws.Range("L136").Value = 9760.5
$ws.Range("M136").Value = -322.1539199999997
$ws.Range("N136").Value = -14860.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1869
$ws.Range("I20").Value = 1851
$ws.Range("J20").Value = 1897.8
$ws.Range("K20").Value = 1851
$ws.Range("L20").Value = 1897.8
$ws.Range("M20").Value = -1604
$ws.Range("N20").Value = -2391.8
$ws.Range("H105").Value = 111113260
$ws.Range("I105").Value = 142859500
$ws.Range("J105").Value = 1455.5
$ws.Range("K105").Value = 142859500
$ws.Range("L105").Value = 1455.5
$ws.Range("M105").Value = -142857753
$ws.Range("N105").Value = -4949.5
$ws.Range("H107").Value = 2075
$ws.Range("J107").Value = 2450
$ws.Range("L107").Value = 2450
$ws.Range("N107").Value = -6290

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 626
$ws.Range("I2").Value = 410
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 410
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = -297
$ws.Range("N2").Value = -1176
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H5").Value = 1735.6666
$ws.Range("I5").Value = 1735.6666
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1735.6666
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1623.6666
$ws.Range("N5").Value = $null
$ws.Range("H10").Value = 663.3333
$ws.Range("I10").Value = 663.3333
$ws.Range("K10").Value = 663.3333
$ws.Range("M10").Value = -524.3333
$ws.Range("H31").Value = 1819.919
$ws.Range("I31").Value = 920.9524
$ws.Range("J31").Value = 2999.8125
$ws.Range("K31").Value = 920.9524
$ws.Range("L31").Value = 2999.8125
$ws.Range("M31").Value = -625.9524
$ws.Range("N31").Value = -3589.8125
$ws.Range("H34").Value = 1819.919
$ws.Range("I34").Value = 920.9524
$ws.Range("J34").Value = 2999.8125
$ws.Range("K34").Value = 920.9524
$ws.Range("L34").Value = 2999.8125
$ws.Range("M34").Value = -718.9524
$ws.Range("N34").Value = -3403.8125
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16316
$ws.Range("H107").Value = 579.62067
$ws.Range("I107").Value = 420.66666
$ws.Range("J107").Value = 839.7273
$ws.Range("K107").Value = 420.66666
$ws.Range("L107").Value = 839.7273
$ws.Range("M107").Value = 1499.33334
$ws.Range("N107").Value = -4679.7273

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 994.5
$ws.Range("I14").Value = 994.5
$ws.Range("K14").Value = 2983.5
$ws.Range("M14").Value = -2810.5
$ws.Range("H131").Value = 22728698
$ws.Range("J131").Value = 1538.55
$ws.Range("L131").Value = 4615.65
$ws.Range("N131").Value = -14695.65
$ws.Range("H139").Value = 2626.1667
$ws.Range("I139").Value = 4034.875
$ws.Range("J139").Value = 1499.2
$ws.Range("K139").Value = 12104.625
$ws.Range("L139").Value = 4497.6
$ws.Range("M139").Value = -6964.625
$ws.Range("N139").Value = -14777.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21434920
$ws.Range("I70").Value = 19236224
$ws.Range("K70").Value = 19236224
$ws.Range("M70").Value = -19235954
$ws.Range("H73").Value = 21434920
$ws.Range("I73").Value = 19236224
$ws.Range("K73").Value = 19236224
$ws.Range("M73").Value = -19235288

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2278.4285
$ws.Range("I7").Value = 2305.5
$ws.Range("K7").Value = 2305.5
$ws.Range("M7").Value = -2193.5
$ws.Range("H40").Value = 7001.25
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 8668.333000000001
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 8668.333000000001
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -8940.333000000001
$ws.Range("H122").Value = 27781306
$ws.Range("I122").Value = 50003400
$ws.Range("J122").Value = 3688.75
$ws.Range("K122").Value = 150010200
$ws.Range("L122").Value = 11066.25
$ws.Range("M122").Value = -150007750
$ws.Range("N122").Value = -15966.25
$ws.Range("H126").Value = 2278.4285
$ws.Range("I126").Value = 2305.5
$ws.Range("K126").Value = 6916.5
$ws.Range("M126").Value = -4446.5
$ws.Range("H136").Value = 2082
$ws.Range("I136").Value = 1817.3334
$ws.Range("K136").Value = 5452.0002
$ws.Range("M136").Value = -2902.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 883.3913
$ws.Range("I136").Value = 725.2941
$ws.Range("J136").Value = 1331.3334
$ws.Range("K136").Value = 2175.8823
$ws.Range("L136").Value = 3994.0002
$ws.Range("M136").Value = 374.1177000000002
$ws.Range("N136").Value = -9094.0002
